$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content edit: patent-notification test case (row 22) gained a third
# sub-case ("...likes his comment an patent" / OPQA-3951) appended to the
# existing pipe-delimited Jira-id / description cells. ---
$ws.Cells.Item(22, 2).Value2 = "OPQA-1433||OPQA-1432||OPQA-3951"
$ws.Cells.Item(22, 3).Value2 = "Verify that user receives a notification when someone he is following comments on a patent||Verify that user receives a notification when someone comments on a patent contained in his watchlist||Verify that user receives a notification if someone likes his comment an patent"

# --- Row heights: these wrap-text rows re-measured slightly (content grew
# on row 22 and the sheet was re-laid-out), matching the recorded cache. ---
$ws.Rows.Item(8).RowHeight = 28.8
$ws.Rows.Item(9).RowHeight = 57.6
$ws.Rows.Item(11).RowHeight = 28.8
$ws.Rows.Item(21).RowHeight = 43.2
$ws.Rows.Item(22).RowHeight = 28.8

# --- View state: selection/scroll moved to the edited cell. ---
$aw = $excel.ActiveWindow
$aw.ScrollRow = 9
$aw.ScrollColumn = 3
$ws.Range("C22").Select() | Out-Null
